$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings
# (e.g. "0.9995", "314.32") are stored as text, matching the
# original inline/shared-string cell type instead of being
# auto-converted to numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.668.04"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.705.46"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "314.32"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "0.9991"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "0.3974"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "0.4061"
$ws.Range("E8").Value = "  +1.64%  "
$ws.Range("D9").Value = "0.9995"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "1.512"
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").Value = "53.27"
$ws.Range("E11").Value = "  +9.32%  "
$ws.Range("D12").Value = "0.08819"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "7.320"
$ws.Range("E13").Value = "  +10.69%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "0.00001328"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "7.580"
$ws.Range("E16").Value = "  +4.85%  "
$ws.Range("D17").Value = "1.700.24"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "101.11"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "0.07121"
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("D20").Value = "19.62"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "6.785"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "0.9994"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "14.18"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").Value = "24.656.90"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "3.003"
$ws.Range("E25").Value = "  +6.90%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "22.54"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "159.97"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "5.136"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").Value = "134.00"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "7.423"
$ws.Range("E31").Value = "  +27.07%  "
$ws.Range("D32").Value = "1.888.53"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("D33").Value = "1.088"
$ws.Range("E33").Value = "  -7.74%  "
$ws.Range("D34").Value = "0.08725"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "7.277"
$ws.Range("E35").Value = "  +17.77%  "
$ws.Range("D36").Value = "11.17"
$ws.Range("E36").Value = "  +2.26%  "
$ws.Range("D37").Value = "1.961"
$ws.Range("E37").Value = "  +6.99%  "
$ws.Range("D38").Value = "0.2737"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").Value = "14.84"
$ws.Range("E39").Value = "  -3.52%  "
$ws.Range("D40").Value = "0.02781"
$ws.Range("E40").Value = "  +9.51%  "
$ws.Range("D41").Value = "0.09023"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "1.483"
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "0.7709"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "0.7209"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "15.66"
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("D46").Value = "2.467"
$ws.Range("E46").Value = "  +1.75%  "
$ws.Range("D47").Value = "4.174"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "0.9990"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "141.57"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "1.302"
$ws.Range("E50").Value = "  +13.27%  "
$ws.Range("D51").Value = "0.00000000373"
$ws.Range("E51").Value = "  +2.61%  "

# Restore the default (Normal) style on column D so no stray
# number-format style lingers on the cells themselves.
$dRange.Style = "Normal"
